# Auto-generated Excel COM-interop script applying numeric updates
# to the Gilgamesh_Profits workbook (per scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 4806.1577  # H40
$ws.Cells.Item(40, 9).Value = 4732  # I40
$ws.Cells.Item(40, 10).Value = 4888.5557  # J40
$ws.Cells.Item(40, 11).Value = 4732  # K40
$ws.Cells.Item(40, 12).Value = 4888.5557  # L40
$ws.Cells.Item(40, 13).Value = -4557  # M40
$ws.Cells.Item(40, 14).Value = -5238.5557  # N40
$ws.Cells.Item(70, 8).Value = 6334  # H70
$ws.Cells.Item(70, 10).Value = 4500  # J70
$ws.Cells.Item(70, 12).Value = 13500  # L70
$ws.Cells.Item(70, 14).Value = -14040  # N70
$ws.Cells.Item(73, 8).Value = 6334  # H73
$ws.Cells.Item(73, 10).Value = 4500  # J73
$ws.Cells.Item(73, 12).Value = 13500  # L73
$ws.Cells.Item(73, 14).Value = -15372  # N73
$ws.Cells.Item(87, 8).Value = 313128.44  # H87
$ws.Cells.Item(87, 10).Value = 313128.44  # J87
$ws.Cells.Item(87, 12).Value = 313128.44  # L87
$ws.Cells.Item(87, 14).Value = -315624.44  # N87
$ws.Cells.Item(90, 8).Value = 313128.44  # H90
$ws.Cells.Item(90, 10).Value = 313128.44  # J90
$ws.Cells.Item(90, 12).Value = 939385.3200000001  # L90
$ws.Cells.Item(90, 14).Value = -951865.3200000001  # N90
$ws.Cells.Item(107, 8).Value = 789.1667  # H107
$ws.Cells.Item(107, 9).Value = 732.41174  # I107
$ws.Cells.Item(107, 10).Value = 927  # J107
$ws.Cells.Item(107, 11).Value = 732.41174  # K107
$ws.Cells.Item(107, 12).Value = 927  # L107
$ws.Cells.Item(107, 13).Value = 1187.58826  # M107
$ws.Cells.Item(107, 14).Value = -4767  # N107
$ws.Cells.Item(115, 8).Value = 889.25  # H115
$ws.Cells.Item(115, 10).Value = 0  # J115
$ws.Cells.Item(115, 12).Value = 0  # L115
$ws.Cells.Item(115, 14).ClearContents()  # N115
$ws.Cells.Item(118, 8).Value = 1225.8462  # H118
$ws.Cells.Item(118, 9).Value = 893.7  # I118
$ws.Cells.Item(118, 11).Value = 2681.1  # K118
$ws.Cells.Item(118, 13).Value = -1024.1  # M118
$ws.Cells.Item(126, 8).Value = 77739.336  # H126
$ws.Cells.Item(126, 10).Value = 77739.336  # J126
$ws.Cells.Item(126, 12).Value = 77739.336  # L126
$ws.Cells.Item(126, 14).Value = -87619.336  # N126
$ws.Cells.Item(132, 8).Value = 4813.3687  # H132
$ws.Cells.Item(132, 9).Value = 4813.3687  # I132
$ws.Cells.Item(132, 11).Value = 14440.1061  # K132
$ws.Cells.Item(132, 13).Value = -11910.1061  # M132
$ws.Cells.Item(138, 8).Value = 443571.06  # H138
$ws.Cells.Item(138, 9).Value = 6497.6665  # I138
$ws.Cells.Item(138, 10).Value = 560123.9399999999  # J138
$ws.Cells.Item(138, 11).Value = 19492.9995  # K138
$ws.Cells.Item(138, 12).Value = 1680371.82  # L138
$ws.Cells.Item(138, 13).Value = -14352.9995  # M138
$ws.Cells.Item(138, 14).Value = -1690651.82  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7673.013  # H32
$ws.Cells.Item(32, 9).Value = 6435.4  # I32
$ws.Cells.Item(32, 11).Value = 6435.4  # K32
$ws.Cells.Item(32, 13).Value = -6148.4  # M32
$ws.Cells.Item(45, 8).Value = 27491.166  # H45
$ws.Cells.Item(45, 10).Value = 5380.5  # J45
$ws.Cells.Item(45, 12).Value = 5380.5  # L45
$ws.Cells.Item(45, 14).Value = -6134.5  # N45
$ws.Cells.Item(61, 8).Value = 5921.0527  # H61
$ws.Cells.Item(61, 9).Value = 5107.143  # I61
$ws.Cells.Item(61, 10).Value = 8200  # J61
$ws.Cells.Item(61, 11).Value = 5107.143  # K61
$ws.Cells.Item(61, 12).Value = 8200  # L61
$ws.Cells.Item(61, 13).Value = -4895.143  # M61
$ws.Cells.Item(61, 14).Value = -8624  # N61
$ws.Cells.Item(74, 8).Value = 296455.94  # H74
$ws.Cells.Item(74, 9).Value = 509019.47  # I74
$ws.Cells.Item(74, 11).Value = 509019.47  # K74
$ws.Cells.Item(74, 13).Value = -508145.47  # M74
$ws.Cells.Item(77, 8).Value = 296455.94  # H77
$ws.Cells.Item(77, 9).Value = 509019.47  # I77
$ws.Cells.Item(77, 11).Value = 2545097.35  # K77
$ws.Cells.Item(77, 13).Value = -2540729.35  # M77
$ws.Cells.Item(97, 8).Value = 1844.4375  # H97
$ws.Cells.Item(97, 9).Value = 2083.5833  # I97
$ws.Cells.Item(97, 10).Value = 1127  # J97
$ws.Cells.Item(97, 11).Value = 2083.5833  # K97
$ws.Cells.Item(97, 12).Value = 1127  # L97
$ws.Cells.Item(97, 13).Value = -1587.5833  # M97
$ws.Cells.Item(97, 14).Value = -2119  # N97
$ws.Cells.Item(102, 8).Value = 4911.6895  # H102
$ws.Cells.Item(102, 9).Value = 4968.1113  # I102
$ws.Cells.Item(102, 10).Value = 4150  # J102
$ws.Cells.Item(102, 11).Value = 4968.1113  # K102
$ws.Cells.Item(102, 12).Value = 4150  # L102
$ws.Cells.Item(102, 13).Value = -3346.1113  # M102
$ws.Cells.Item(102, 14).Value = -7394  # N102
$ws.Cells.Item(136, 8).Value = 5921.0527  # H136
$ws.Cells.Item(136, 9).Value = 5107.143  # I136
$ws.Cells.Item(136, 10).Value = 8200  # J136
$ws.Cells.Item(136, 11).Value = 15321.429  # K136
$ws.Cells.Item(136, 12).Value = 24600  # L136
$ws.Cells.Item(136, 13).Value = -12771.429  # M136
$ws.Cells.Item(136, 14).Value = -29700  # N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 26885364  # H20
$ws.Cells.Item(20, 9).Value = 33337420  # I20
$ws.Cells.Item(20, 10).Value = 1796.5  # J20
$ws.Cells.Item(20, 11).Value = 33337420  # K20
$ws.Cells.Item(20, 12).Value = 1796.5  # L20
$ws.Cells.Item(20, 13).Value = -33337173  # M20
$ws.Cells.Item(20, 14).Value = -2290.5  # N20
$ws.Cells.Item(31, 8).Value = 4474.3164  # H31
$ws.Cells.Item(31, 9).Value = 3741.4897  # I31
$ws.Cells.Item(31, 10).Value = 5671.2666  # J31
$ws.Cells.Item(31, 11).Value = 3741.4897  # K31
$ws.Cells.Item(31, 12).Value = 5671.2666  # L31
$ws.Cells.Item(31, 13).Value = -3446.4897  # M31
$ws.Cells.Item(31, 14).Value = -6261.2666  # N31
$ws.Cells.Item(34, 8).Value = 4474.3164  # H34
$ws.Cells.Item(34, 9).Value = 3741.4897  # I34
$ws.Cells.Item(34, 10).Value = 5671.2666  # J34
$ws.Cells.Item(34, 11).Value = 3741.4897  # K34
$ws.Cells.Item(34, 12).Value = 5671.2666  # L34
$ws.Cells.Item(34, 13).Value = -3539.4897  # M34
$ws.Cells.Item(34, 14).Value = -6075.2666  # N34
$ws.Cells.Item(94, 8).Value = 2340.4  # H94
$ws.Cells.Item(94, 9).Value = 1352  # I94
$ws.Cells.Item(94, 11).Value = 1352  # K94
$ws.Cells.Item(94, 13).Value = -901  # M94
$ws.Cells.Item(99, 8).Value = 5918.15  # H99
$ws.Cells.Item(99, 9).Value = 6120.625  # I99
$ws.Cells.Item(99, 11).Value = 6120.625  # K99
$ws.Cells.Item(99, 13).Value = -4622.625  # M99
$ws.Cells.Item(107, 8).Value = 1006.75  # H107
$ws.Cells.Item(107, 9).Value = 1075.3334  # I107
$ws.Cells.Item(107, 10).Value = 801  # J107
$ws.Cells.Item(107, 11).Value = 1075.3334  # K107
$ws.Cells.Item(107, 12).Value = 801  # L107
$ws.Cells.Item(107, 13).Value = 844.6666  # M107
$ws.Cells.Item(107, 14).Value = -4641  # N107
$ws.Cells.Item(126, 8).Value = 5918.15  # H126
$ws.Cells.Item(126, 9).Value = 6120.625  # I126
$ws.Cells.Item(126, 11).Value = 18361.875  # K126
$ws.Cells.Item(126, 13).Value = -15891.875  # M126
$ws.Cells.Item(132, 8).Value = 6668771.5  # H132
$ws.Cells.Item(132, 9).Value = 7464705.5  # I132
$ws.Cells.Item(132, 10).Value = 2826.875  # J132
$ws.Cells.Item(132, 11).Value = 22394116.5  # K132
$ws.Cells.Item(132, 12).Value = 8480.625  # L132
$ws.Cells.Item(132, 13).Value = -22391586.5  # M132
$ws.Cells.Item(132, 14).Value = -13540.625  # N132
$ws.Cells.Item(134, 8).Value = 3648.7144  # H134
$ws.Cells.Item(134, 9).Value = 2196.3333  # I134
$ws.Cells.Item(134, 11).Value = 6588.999899999999  # K134
$ws.Cells.Item(134, 13).Value = -4053.999899999999  # M134
$ws.Cells.Item(141, 8).Value = 284503.47  # H141
$ws.Cells.Item(141, 10).Value = 284503.47  # J141
$ws.Cells.Item(141, 12).Value = 284503.47  # L141
$ws.Cells.Item(141, 14).Value = -294863.47  # N141

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 150  # H7
$ws.Cells.Item(7, 9).Value = 150  # I7
$ws.Cells.Item(7, 11).Value = 450  # K7
$ws.Cells.Item(7, 13).Value = -338  # M7
$ws.Cells.Item(39, 8).Value = 11176.615  # H39
$ws.Cells.Item(39, 10).Value = 11176.615  # J39
$ws.Cells.Item(39, 12).Value = 33529.845  # L39
$ws.Cells.Item(39, 14).Value = -34117.845  # N39
$ws.Cells.Item(55, 8).Value = 4065.5789  # H55
$ws.Cells.Item(55, 10).Value = 4553.3125  # J55
$ws.Cells.Item(55, 12).Value = 13659.9375  # L55
$ws.Cells.Item(55, 14).Value = -14013.9375  # N55
$ws.Cells.Item(82, 8).Value = 11400  # H82
$ws.Cells.Item(82, 10).Value = 11400  # J82
$ws.Cells.Item(82, 12).Value = 34200  # L82
$ws.Cells.Item(82, 14).Value = -35012  # N82
$ws.Cells.Item(85, 8).Value = 11400  # H85
$ws.Cells.Item(85, 10).Value = 11400  # J85
$ws.Cells.Item(85, 12).Value = 34200  # L85
$ws.Cells.Item(85, 14).Value = -37008  # N85
$ws.Cells.Item(92, 8).Value = 592.25  # H92
$ws.Cells.Item(92, 9).Value = 593.3333  # I92
$ws.Cells.Item(92, 10).Value = 589  # J92
$ws.Cells.Item(92, 11).Value = 1779.9999  # K92
$ws.Cells.Item(92, 12).Value = 1767  # L92
$ws.Cells.Item(92, 13).Value = -531.9999  # M92
$ws.Cells.Item(92, 14).Value = -4263  # N92

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 142860750  # H80
$ws.Cells.Item(80, 10).Value = 4251.4  # J80
$ws.Cells.Item(80, 12).Value = 4251.4  # L80
$ws.Cells.Item(80, 14).Value = -6247.4  # N80
$ws.Cells.Item(83, 8).Value = 142860750  # H83
$ws.Cells.Item(83, 10).Value = 4251.4  # J83
$ws.Cells.Item(83, 12).Value = 21257  # L83
$ws.Cells.Item(83, 14).Value = -31241  # N83
$ws.Cells.Item(93, 8).Value = 89998.5  # H93
$ws.Cells.Item(93, 10).Value = 89998.5  # J93
$ws.Cells.Item(93, 12).Value = 89998.5  # L93
$ws.Cells.Item(93, 14).Value = -93742.5  # N93
$ws.Cells.Item(122, 8).Value = 2126.2  # H122
$ws.Cells.Item(122, 9).Value = 2293.25  # I122
$ws.Cells.Item(122, 10).Value = 1875.625  # J122
$ws.Cells.Item(122, 11).Value = 6879.75  # K122
$ws.Cells.Item(122, 12).Value = 5626.875  # L122
$ws.Cells.Item(122, 13).Value = -4429.75  # M122
$ws.Cells.Item(122, 14).Value = -10526.875  # N122

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 2330.6667  # H93
$ws.Cells.Item(93, 9).Value = 2501  # I93
$ws.Cells.Item(93, 11).Value = 2501  # K93
$ws.Cells.Item(93, 13).Value = -1253  # M93
$ws.Cells.Item(132, 8).Value = 2947.4  # H132
$ws.Cells.Item(132, 9).Value = 2477.9614  # I132
$ws.Cells.Item(132, 10).Value = 5998.75  # J132
$ws.Cells.Item(132, 11).Value = 7433.8842  # K132
$ws.Cells.Item(132, 12).Value = 17996.25  # L132
$ws.Cells.Item(132, 13).Value = -4903.8842  # M132
$ws.Cells.Item(132, 14).Value = -23056.25  # N132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 5797.4  # H62
$ws.Cells.Item(62, 9).Value = 5797.4  # I62
$ws.Cells.Item(62, 10).Value = 0  # J62
$ws.Cells.Item(62, 11).Value = 5797.4  # K62
$ws.Cells.Item(62, 12).Value = 0  # L62
$ws.Cells.Item(62, 13).Value = -5173.4  # M62
$ws.Cells.Item(62, 14).ClearContents()  # N62
$ws.Cells.Item(65, 8).Value = 5797.4  # H65
$ws.Cells.Item(65, 9).Value = 5797.4  # I65
$ws.Cells.Item(65, 10).Value = 0  # J65
$ws.Cells.Item(65, 11).Value = 28987  # K65
$ws.Cells.Item(65, 12).Value = 0  # L65
$ws.Cells.Item(65, 13).Value = -25867  # M65
$ws.Cells.Item(65, 14).ClearContents()  # N65
$ws.Cells.Item(81, 8).Value = 4734.4116  # H81
$ws.Cells.Item(81, 9).Value = 3407.7273  # I81
$ws.Cells.Item(81, 10).Value = 7166.6665  # J81
$ws.Cells.Item(81, 11).Value = 6815.4546  # K81
$ws.Cells.Item(81, 12).Value = 14333.333  # L81
$ws.Cells.Item(81, 13).Value = -5754.4546  # M81
$ws.Cells.Item(81, 14).Value = -16455.333  # N81
$ws.Cells.Item(84, 8).Value = 4734.4116  # H84
$ws.Cells.Item(84, 9).Value = 3407.7273  # I84
$ws.Cells.Item(84, 10).Value = 7166.6665  # J84
$ws.Cells.Item(84, 11).Value = 34077.273  # K84
$ws.Cells.Item(84, 12).Value = 71666.66500000001  # L84
$ws.Cells.Item(84, 13).Value = -28773.273  # M84
$ws.Cells.Item(84, 14).Value = -82274.66500000001  # N84
$ws.Cells.Item(122, 8).Value = 12503267  # H122
$ws.Cells.Item(122, 9).Value = 2896.2307  # I122
$ws.Cells.Item(122, 10).Value = 35718244  # J122
$ws.Cells.Item(122, 11).Value = 8688.6921  # K122
$ws.Cells.Item(122, 12).Value = 107154732  # L122
$ws.Cells.Item(122, 13).Value = -6238.6921  # M122
$ws.Cells.Item(122, 14).Value = -107159632  # N122
$ws.Cells.Item(126, 8).Value = 2116.5417  # H126
$ws.Cells.Item(126, 9).Value = 2116.5417  # I126
$ws.Cells.Item(126, 11).Value = 6349.625100000001  # K126
$ws.Cells.Item(126, 13).Value = -3879.625100000001  # M126
$ws.Cells.Item(133, 8).Value = 81698.75  # H133
$ws.Cells.Item(133, 10).Value = 87748  # J133
$ws.Cells.Item(133, 12).Value = 87748  # L133
$ws.Cells.Item(133, 14).Value = -97868  # N133
$ws.Cells.Item(136, 8).Value = 24391840  # H136
$ws.Cells.Item(136, 9).Value = 27028434  # I136
$ws.Cells.Item(136, 10).Value = 3349.75  # J136
$ws.Cells.Item(136, 11).Value = 81085302  # K136
$ws.Cells.Item(136, 12).Value = 10049.25  # L136
$ws.Cells.Item(136, 13).Value = -81082752  # M136
$ws.Cells.Item(136, 14).Value = -15149.25  # N136
